# Build out the "current movement" / "closing" columns (H:K) of the trial
# table and fill in the figures that go with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cells on row 1 (H1:K1)
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "current movement debt"
$ws.Range("I1").Value = "current movement credit"
$ws.Range("J1").Value = "closing debt"
$ws.Range("K1").Value = "closing credit"

# Clone the look of the existing header cells (bold white text, centered,
# dark fill) onto the new header cells by copying A1's format.
$ws.Range("A1").Copy()
$ws.Range("H1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# I1 is made of two runs in the source workbook ("current movement" + "
# credit"); reproduce that as rich text, both runs keeping the bold/white
# header styling.
$run1 = $ws.Range("I1").Characters(1, 16)
$run1.Font.Bold = $true
$run1.Font.Color = 16777215
$run1.Font.Name = "Calibri"
$run1.Font.Size = 11

$run2 = $ws.Range("I1").Characters(17, 7)
$run2.Font.Bold = $true
$run2.Font.Color = 16777215
$run2.Font.Name = "Calibri"
$run2.Font.Size = 11

# ---------------------------------------------------------------------
# 2. Figures for the new columns
# ---------------------------------------------------------------------
$ws.Range("J6").Value = 15000
$ws.Range("K7").Value = 6000
$ws.Range("I11").Value = 5000
$ws.Range("K11").Value = 55000
$ws.Range("F21").Value = 4500
$ws.Range("J21").Value = 1000
$ws.Range("G24").Value = 9500
$ws.Range("I24").Value = 6000
$ws.Range("J24").Value = 500
$ws.Range("H25").Value = 8000
$ws.Range("G27").Value = 8500
$ws.Range("J27").Value = 6000
$ws.Range("F28").Value = 1500
$ws.Range("K28").Value = 2300

# ---------------------------------------------------------------------
# 3. Column widths for the new columns
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 21.333333333333332
$ws.Columns.Item(9).ColumnWidth = 22.333333333333332
$ws.Columns.Item(10).ColumnWidth = 12.166666666666666
$ws.Columns.Item(11).ColumnWidth = 13.166666666666666
$ws.Columns.Item(12).ColumnWidth = 10.666666666666666

# ---------------------------------------------------------------------
# 4. Selection / active cell to match the final view
# ---------------------------------------------------------------------
$ws.Range("K16").Select()

Write-Host "trial table columns H:K populated"
